$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 16  # row A=2 "北京·2024爬宠狂欢节首展·森临奇遇记": 15 -> 16
$ws.Range("F4").Value = 1348  # row A=3 "北京·第三届次元潮流动漫嘉年华": 1347 -> 1348
$ws.Range("F6").Value = 7698  # row A=5 "北京·Aw动漫游戏嘉年华7th- 来自异世界的召唤": 7699 -> 7698
$ws.Range("F9").Value = 2101  # row A=8 "北京·TCS卡牌嘉年华": 2100 -> 2101
$ws.Range("F10").Value = 8490  # row A=9 "北京·thebONE✖️GOJO超次元嘉年华02": 8487 -> 8490
$ws.Range("F13").Value = 77  # row A=12 "北京·THP 03 白兔茶话会": 74 -> 77
$ws.Range("F14").Value = 5687  # row A=13 "北京·thebONE×Ilike动漫游戏嘉年华S4": 5684 -> 5687
$ws.Range("F16").Value = 2648  # row A=15 "北京·万游引力国潮动漫嘉年华s6": 2646 -> 2648
$ws.Range("F17").Value = 1151  # row A=16 "北京·排球少年ONLY": 1149 -> 1151
$ws.Range("F19").Value = 350  # row A=18 "北京·首届明日方舟only展·春和京明": 349 -> 350
$ws.Range("F22").Value = 35  # row A=21 "北京·2024图书市集春季场": 34 -> 35
$ws.Range("F24").Value = 3630  # row A=23 "北京·IDOx梦次元动漫游戏嘉年华3rd": 3621 -> 3630
$ws.Range("F29").Value = 3123  # row A=28 "北京·第15届IJOY漫展xCGF游戏节": 3118 -> 3123
$ws.Range("F30").Value = 57  # row A=29 "北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 小N&小敢 专场活动": 56 -> 57
$ws.Range("F31").Value = 197  # row A=30 "北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 杨天翔 专场活动": 190 -> 197
$ws.Range("F32").Value = 357  # row A=31 "北京·Yok运动番Only": 355 -> 357
$ws.Range("F34").Value = 329  # row A=33 "北京·第五人格Only": 327 -> 329
$ws.Range("F35").Value = 783  # row A=34 "北京·第22届中国国际模型博览会": 738 -> 783
$ws.Range("F36").Value = 671  # row A=35 "北京·QMQ动漫游戏嘉年华": 670 -> 671
$ws.Range("F39").Value = 2211  # row A=38 "北京·IDO动漫游戏嘉年华45th": 2136 -> 2211
$ws.Range("F40").Value = 51  # row A=39 "北京·IDO动漫游戏嘉年华45th同人创作大会": 50 -> 51
$ws.Range("F43").Value = 3057  # row A=42 "北京·第16届IJOY漫展XCGF游戏节": 3047 -> 3057
$ws.Range("F45").Value = 2297  # row A=44 "北京·原神only3.0": 2296 -> 2297
$ws.Range("F47").Value = 32  # row A=46 "北京·原神only3.0——蛋黄mayo签售会": 31 -> 32

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 132  # row A=2 "北京·「京都动画X春日计划」漫展感管弦室内乐音乐会 ": 131 -> 132
$ws.Range("F10").Value = 2  # row A=9 "北京·《国风大赏》大型国潮音乐会×郑州歌舞剧院《唐宫夜宴》": 1 -> 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1334  # row A=2 "北京·盗墓笔记官方授权主题店": 1333 -> 1334

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1334  # row A=2 "北京·盗墓笔记官方授权主题店": 1333 -> 1334
$ws.Range("F4").Value = 16  # row A=3 "北京·2024爬宠狂欢节首展·森临奇遇记": 15 -> 16
$ws.Range("F5").Value = 1348  # row A=4 "北京·第三届次元潮流动漫嘉年华": 1347 -> 1348
$ws.Range("F6").Value = 7698  # row A=5 "北京·Aw动漫游戏嘉年华7th- 来自异世界的召唤": 7699 -> 7698
$ws.Range("F9").Value = 2101  # row A=8 "北京·TCS卡牌嘉年华": 2100 -> 2101
$ws.Range("F10").Value = 8490  # row A=9 "北京·thebONE✖️GOJO超次元嘉年华02": 8487 -> 8490
$ws.Range("F12").Value = 77  # row A=11 "北京·THP 03 白兔茶话会": 74 -> 77
$ws.Range("F13").Value = 5687  # row A=12 "北京·thebONE×Ilike动漫游戏嘉年华S4": 5684 -> 5687
$ws.Range("F15").Value = 2648  # row A=14 "北京·万游引力国潮动漫嘉年华s6": 2646 -> 2648
$ws.Range("F16").Value = 1151  # row A=15 "北京·排球少年ONLY": 1149 -> 1151
$ws.Range("F21").Value = 35  # row A=20 "北京·2024图书市集春季场": 34 -> 35
$ws.Range("F22").Value = 132  # row A=21 "北京·「京都动画X春日计划」漫展感管弦室内乐音乐会 ": 131 -> 132
$ws.Range("F25").Value = 3630  # row A=24 "北京·IDOx梦次元动漫游戏嘉年华3rd": 3621 -> 3630
$ws.Range("F30").Value = 3123  # row A=29 "北京·第15届IJOY漫展xCGF游戏节": 3118 -> 3123
$ws.Range("F31").Value = 57  # row A=30 "北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 小N&小敢 专场活动": 56 -> 57
$ws.Range("F32").Value = 357  # row A=31 "北京·Yok运动番Only": 355 -> 357
$ws.Range("F34").Value = 329  # row A=33 "北京·第五人格Only": 327 -> 329
$ws.Range("F36").Value = 783  # row A=35 "北京·第22届中国国际模型博览会": 738 -> 783
$ws.Range("F37").Value = 671  # row A=36 "北京·QMQ动漫游戏嘉年华": 670 -> 671
$ws.Range("F41").Value = 2212  # row A=40 "北京·IDO动漫游戏嘉年华45th": 2137 -> 2212
$ws.Range("F42").Value = 51  # row A=41 "北京·IDO动漫游戏嘉年华45th同人创作大会": 50 -> 51
$ws.Range("F45").Value = 3057  # row A=44 "北京·第16届IJOY漫展XCGF游戏节": 3047 -> 3057
$ws.Range("F46").Value = 2297  # row A=45 "北京·原神only3.0": 2296 -> 2297
$ws.Range("F49").Value = 2  # row A=48 "北京·《国风大赏》大型国潮音乐会×郑州歌舞剧院《唐宫夜宴》": 1 -> 2
